$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.186.10"
$ws.Range("D3").Value = "2.235.42"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.44"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.601"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0957"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.52%  "
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").Value = "2.573.72"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.838"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.52%  "
$ws.Range("D17").Value = "2.227.30"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "42.105.09"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +13.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "230.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.18%  "
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.04%  "
$ws.Range("E29").Value = "  -2.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("E32").Value = "  -4.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0804"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("E36").Value = "  -7.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0304"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.49%  "
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.36%  "
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").Value = "2.444.07"
$ws.Range("E51").Value = "  -0.26%  "
